$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Remove the standalone "Meta description: ..." paragraph that currently
#    sits right after the H1 title.
# ---------------------------------------------------------------------------
$metaPara = $d.Paragraphs.Item(2)
if ($metaPara.Range.Text -match "Meta description") {
    $metaPara.Range.Delete()
}

# ---------------------------------------------------------------------------
# 2. Insert a new bold paragraph ("Play Bugs Money for Free - Exciting
#    Features and Impeccable Graphics") right before the final paragraph
#    (the one that held the image-generation prompt, in italics).
# ---------------------------------------------------------------------------
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($count)
$insertPoint = $d.Range($lastPara.Range.Start, $lastPara.Range.Start)

$xmlSnippet = '<?xml version="1.0"?><?mso-application progid="Word.Document"?>' +
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData>' +
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:body>' +
    '<w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Bugs Money for Free - Exciting Features and Impeccable Graphics</w:t></w:r></w:p>' +
    '<w:p/>' +
    '</w:body></w:document>' +
    '</pkg:xmlData></pkg:part></pkg:package>'

$insertPoint.InsertXML($xmlSnippet) | Out-Null

# InsertXML splits the destination paragraph into two; the first one carries
# our new bold text, the second is an empty spacer paragraph that needs to
# be dropped so the original final paragraph directly follows our insert.
$spacerPara = $d.Paragraphs.Item($count + 1)
$spacerPara.Range.Delete()

# ---------------------------------------------------------------------------
# 3. Replace the text of the last (italic) paragraph with the review blurb,
#    keeping its existing run formatting (italic) untouched.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Create a cartoon-style feature image for Bugs Money that prominently features a happy Maya warrior wearing glasses. The warrior should be standing in a bright green lawn with flowers and bugs surrounding them, holding a honeycomb with money flying out of it in the background. It should convey the fun and colorful nature of the game and showcase the potential for big wins. The image should be eye-catching and playful, enticing players to try their luck with Bugs Money.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Read our review of Bugs Money and play for free. Exciting features and impeccable graphics, including Glow Wilds and Free Spins.",
    2) | Out-Null
